$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after the "Contact" row (row 10), pushing existing rows down.
$ws.Range("A11:B11").Insert()

# Copy formatting from the row above (Contact, row 10) to keep the look consistent,
# restricted to the used columns only.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with Jurisdiction / iso:code:3166:FR
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Update Version value (row 3)
$ws.Range("B3").Value = "0.2.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"
